$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table data (player name, position, team) in the new row order.
# Column A (player names) happens to end up in the same literal order as
# the original rows, while columns B (Pozisyon) and C (Takım) are
# re-associated per the new row layout.
$data = @(
    @("Jamal Murray",        "PG,SG",    "Denver Nuggets"),
    @("Donovan Mitchell",    "PG,SG",    "Cleveland Cavaliers"),
    @("Payton Pritchard",    "PG",       "Boston Celtics"),
    @("Malik Beasley",       "SG,SF",    "Detroit Pistons"),
    @("Dyson Daniels",       "PG,SG,SF", "Atlanta Hawks"),
    @("De'Andre Hunter",     "SF,PF",    "Atlanta Hawks"),
    @("Michael Porter Jr.",  "SF,PF",    "Denver Nuggets"),
    @("Josh Hart",           "SG,SF,PF", "New York Knicks"),
    @("Myles Turner",        "C",        "Indiana Pacers"),
    @("Deandre Ayton",       "C",        "Portland Trail Blazers"),
    @("Victor Wembanyama",   "C",        "San Antonio Spurs"),
    @("Kristaps Porzingis",  "PF,C",     "Boston Celtics"),
    @("Obi Toppin",          "PF",       "Indiana Pacers"),
    @("Domantas Sabonis",    "C",        "Sacramento Kings"),
    @("Bradley Beal",        "PG,SG,SF", "Phoenix Suns"),
    @("Cam Thomas",          "SG,SF",    "Brooklyn Nets"),
    @("Tari Eason",          "SF,PF",    "Houston Rockets")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
